# Data dictionary update: mark a few fields as candidates for removal, and
# rename/clarify the "AU All" uniqueness note.
#
# Commit message: "noted fields that could be removed in data dictionary"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A "remove" markers -------------------------------------------
# Flag the "Fund- # & Name", "AU All" and "budgetType" fields as removable.
$ws.Range("A6").Value  = "remove"
$ws.Range("A17").Value = "remove"
$ws.Range("A18").Value = "remove"

# --- Clarify the uniqueness note for "AU All" / "Year" -------------------
$ws.Range("G17").Value = "AU x year unique"
$ws.Range("G20").Value = "AU x year unique (AU = Service# - Acct#)"

# --- Restore the view/selection state left behind by the edit ------------
$ws.Range("G21").Select()
